$d = $word.ActiveDocument

# 1. Fix "timestap" -> "timestamp" typo (also removes the proofErr spellcheck markers)
$d.Content.Find.Execute("timestap", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "timestamp", 2) | Out-Null

Write-Output "done"
